$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '27.133.70'
Set-TextValue 'E2' '  +0.66%  '
Set-TextValue 'D3' '1.826.16'
Set-TextValue 'E3' '  +0.45%  '
Set-TextValue 'E4' '  +0.88%  '
Set-TextValue 'D5' '312.89'
Set-TextValue 'E5' '  +0.92%  '
Set-TextValue 'E6' '  +0.79%  '
Set-TextValue 'E7' '  +0.10%  '
Set-TextValue 'E8' '  -0.48%  '
Set-TextValue 'D9' '0.07397'
Set-TextValue 'E9' '  +0.60%  '
Set-TextValue 'D10' '0.8795'
Set-TextValue 'E10' '  +0.76%  '
Set-TextValue 'D11' '20.34'
Set-TextValue 'E11' '  +0.30%  '
Set-TextValue 'D12' '1.926.10'
Set-TextValue 'E12' '  +6.12%  '
Set-TextValue 'E13' '  +3.12%  '
Set-TextValue 'D14' '5.380'
Set-TextValue 'E14' '  -0.53%  '
Set-TextValue 'D15' '93.20'
Set-TextValue 'E15' '  +1.96%  '
Set-TextValue 'D16' '6.513'
Set-TextValue 'E16' '  +0.04%  '
Set-TextValue 'D17' '1.008'
Set-TextValue 'E17' '  +0.51%  '
Set-TextValue 'D18' '0.000008701'
Set-TextValue 'E18' '  -0.05%  '
Set-TextValue 'D19' '1.010'
Set-TextValue 'E19' '  +0.81%  '
Set-TextValue 'D20' '27.680.55'
Set-TextValue 'E20' '  +2.64%  '
Set-TextValue 'D21' '14.64'
Set-TextValue 'E21' '  -0.14%  '
Set-TextValue 'D22' '5.240'
Set-TextValue 'E22' '  -1.02%  '
Set-TextValue 'D23' '10.57'
Set-TextValue 'E23' '  -0.27%  '
Set-TextValue 'D24' '2.099.83'
Set-TextValue 'E24' '  +2.97%  '
Set-TextValue 'E25' '  -0.62%  '
Set-TextValue 'D26' '151.61'
Set-TextValue 'E26' '  +0.52%  '
Set-TextValue 'E27' '  +0.67%  '
Set-TextValue 'D28' '2.133'
Set-TextValue 'E28' '  -0.72%  '
Set-TextValue 'D29' '5.174'
Set-TextValue 'E29' '  -1.61%  '
Set-TextValue 'D30' '116.20'
Set-TextValue 'E30' '  -0.69%  '
Set-TextValue 'E31' '  +0.36%  '
Set-TextValue 'E32' '  +0.48%  '
Set-TextValue 'D33' '0.7404'
Set-TextValue 'E33' '  -2.63%  '
Set-TextValue 'D34' '4.501'
Set-TextValue 'E34' '  -0.02%  '
Set-TextValue 'D35' '2.940'
Set-TextValue 'E35' '  +1.04%  '
Set-TextValue 'E36' '  +0.82%  '
Set-TextValue 'D37' '1.088'
Set-TextValue 'E37' '  -0.59%  '
Set-TextValue 'E38' '  +0.02%  '
Set-TextValue 'D39' '0.01942'
Set-TextValue 'E39' '  -0.18%  '
Set-TextValue 'D40' '2.406'
Set-TextValue 'E40' '  +0.33%  '
Set-TextValue 'D41' '2.939'
Set-TextValue 'E41' '  -1.03%  '
Set-TextValue 'D42' '7.183'
Set-TextValue 'E42' '  +0.44%  '
Set-TextValue 'D43' '0.5255'
Set-TextValue 'E43' '  -0.67%  '
Set-TextValue 'D44' '0.1641'
Set-TextValue 'E44' '  -0.79%  '
Set-TextValue 'D45' '8.363'
Set-TextValue 'E45' '  -0.96%  '
Set-TextValue 'D46' '0.4869'
Set-TextValue 'E46' '  -0.04%  '
Set-TextValue 'E47' '  -0.49%  '
Set-TextValue 'E48' '  +0.86%  '
Set-TextValue 'D49' '104.26'
Set-TextValue 'E49' '  +0.76%  '
Set-TextValue 'D50' '1.650'
Set-TextValue 'E50' '  -0.83%  '
Set-TextValue 'D51' '0.06297'
Set-TextValue 'E51' '  +0.00%  '
